# Apply scheduled-runner price/profit updates across the Sheets workbook.
# Mirrors a refreshed market-board data pull: currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and the derived LeveProfit(NQ/HQ) columns (H:N) are
# rewritten per affected leve row; everything else is untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Range("H28").Value = 804.3333
$ws.Range("I28").Value = 915.3125
$ws.Range("J28").Value = 449.2
$ws.Range("K28").Value = 915.3125
$ws.Range("L28").Value = 449.2
$ws.Range("M28").Value = -430.3125
$ws.Range("N28").Value = -1419.2

# row 55
$ws.Range("H55").Value = 149.9375
$ws.Range("I55").Value = 121.375
$ws.Range("J55").Value = 178.5
$ws.Range("K55").Value = 121.375
$ws.Range("L55").Value = 178.5
$ws.Range("M55").Value = 92.625
$ws.Range("N55").Value = -606.5

# row 82
$ws.Range("H82").Value = 18614.143
$ws.Range("J82").Value = 25399.8
$ws.Range("L82").Value = 76199.39999999999
$ws.Range("N82").Value = -77011.39999999999

# row 85
$ws.Range("H85").Value = 18614.143
$ws.Range("J85").Value = 25399.8
$ws.Range("L85").Value = 76199.39999999999
$ws.Range("N85").Value = -79007.39999999999

# row 92
$ws.Range("H92").Value = 7354734.5
$ws.Range("I92").Value = 8065808.5
$ws.Range("K92").Value = 8065808.5
$ws.Range("M92").Value = -8064560.5

# row 111
$ws.Range("H111").Value = 3487.6667
$ws.Range("I111").Value = 1912.7142
$ws.Range("J111").Value = 9000
$ws.Range("K111").Value = 5738.142599999999
$ws.Range("L111").Value = 27000
$ws.Range("M111").Value = -2671.142599999999
$ws.Range("N111").Value = -33134

# row 113
$ws.Range("H113").Value = 3707.889
$ws.Range("I113").Value = 2354.0908
$ws.Range("J113").Value = 5835.2856
$ws.Range("K113").Value = 2354.0908
$ws.Range("L113").Value = 5835.2856
$ws.Range("M113").Value = 899.9092000000001
$ws.Range("N113").Value = -12343.2856

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1814
$ws.Range("I2").Value = 1208.5834
$ws.Range("J2").Value = 4720
$ws.Range("K2").Value = 1208.5834
$ws.Range("L2").Value = 4720
$ws.Range("M2").Value = -1095.5834
$ws.Range("N2").Value = -4946

# row 44
$ws.Range("H44").Value = 34000
$ws.Range("J44").Value = 34000
$ws.Range("L44").Value = 34000
$ws.Range("N44").Value = -34976

# row 45
$ws.Range("H45").Value = 4113.7
$ws.Range("I45").Value = 3242.3
$ws.Range("K45").Value = 3242.3
$ws.Range("M45").Value = -2865.3

# row 55
$ws.Range("H55").Value = 14219.333
$ws.Range("J55").Value = 14219.333
$ws.Range("L55").Value = 14219.333
$ws.Range("N55").Value = -14849.333

# row 74
$ws.Range("H74").Value = 1101
$ws.Range("I74").Value = 1153.44
$ws.Range("J74").Value = 955.3333
$ws.Range("K74").Value = 1153.44
$ws.Range("L74").Value = 955.3333
$ws.Range("M74").Value = -279.4400000000001
$ws.Range("N74").Value = -2703.3333

# row 77
$ws.Range("H77").Value = 1101
$ws.Range("I77").Value = 1153.44
$ws.Range("J77").Value = 955.3333
$ws.Range("K77").Value = 5767.200000000001
$ws.Range("L77").Value = 4776.6665
$ws.Range("M77").Value = -1399.200000000001
$ws.Range("N77").Value = -13512.6665

# row 80
$ws.Range("H80").Value = 19914.2
$ws.Range("J80").Value = 19914.2
$ws.Range("L80").Value = 19914.2
$ws.Range("N80").Value = -21910.2

# row 83
$ws.Range("H83").Value = 19914.2
$ws.Range("J83").Value = 19914.2
$ws.Range("L83").Value = 59742.60000000001
$ws.Range("N83").Value = -69726.60000000001

# row 97
$ws.Range("H97").Value = 1176.25
$ws.Range("I97").Value = 1001.4286
$ws.Range("K97").Value = 1001.4286
$ws.Range("M97").Value = -505.4286

# row 110
$ws.Range("H110").Value = 10117.333
$ws.Range("I110").Value = 10722.286
$ws.Range("J110").Value = 8000
$ws.Range("K110").Value = 10722.286
$ws.Range("L110").Value = 8000
$ws.Range("M110").Value = -8677.286
$ws.Range("N110").Value = -12090

# row 116
$ws.Range("H116").Value = 1814
$ws.Range("I116").Value = 1208.5834
$ws.Range("J116").Value = 4720
$ws.Range("K116").Value = 1208.5834
$ws.Range("L116").Value = 4720
$ws.Range("M116").Value = 1085.4166
$ws.Range("N116").Value = -9308

# row 122
$ws.Range("H122").Value = 1664.6
$ws.Range("I122").Value = 1749.4286
$ws.Range("J122").Value = 1466.6666
$ws.Range("K122").Value = 5248.2858
$ws.Range("L122").Value = 4399.9998
$ws.Range("M122").Value = -2798.2858
$ws.Range("N122").Value = -9299.9998

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1814
$ws.Range("I3").Value = 1208.5834
$ws.Range("J3").Value = 4720
$ws.Range("K3").Value = 1208.5834
$ws.Range("L3").Value = 4720
$ws.Range("M3").Value = -1094.5834
$ws.Range("N3").Value = -4948

# row 35
$ws.Range("H35").Value = 33388
$ws.Range("J35").Value = 33388
$ws.Range("L35").Value = 33388
$ws.Range("N35").Value = -34008

# row 76
$ws.Range("H76").Value = 9900
$ws.Range("J76").Value = 9900
$ws.Range("L76").Value = 9900
$ws.Range("N76").Value = -10530

# row 79
$ws.Range("H79").Value = 9900
$ws.Range("J79").Value = 9900
$ws.Range("L79").Value = 9900
$ws.Range("N79").Value = -12084

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2467.1943
$ws.Range("I31").Value = 1659.4706
$ws.Range("K31").Value = 1659.4706
$ws.Range("M31").Value = -1364.4706

# row 34
$ws.Range("H34").Value = 2467.1943
$ws.Range("I34").Value = 1659.4706
$ws.Range("K34").Value = 1659.4706
$ws.Range("M34").Value = -1457.4706

# row 60
$ws.Range("H60").Value = 18094.334
$ws.Range("J60").Value = 21094.6
$ws.Range("L60").Value = 21094.6
$ws.Range("N60").Value = -22116.6

# row 68
$ws.Range("H68").Value = 17549.5
$ws.Range("J68").Value = 17549.5
$ws.Range("L68").Value = 17549.5
$ws.Range("N68").Value = -19047.5

# row 71
$ws.Range("H71").Value = 17549.5
$ws.Range("J71").Value = 17549.5
$ws.Range("L71").Value = 52648.5
$ws.Range("N71").Value = -60136.5

# row 109
$ws.Range("H109").Value = 10900
$ws.Range("J109").Value = 10900
$ws.Range("L109").Value = 10900
$ws.Range("N109").Value = -12980

# row 134
$ws.Range("H134").Value = 6086.5713
$ws.Range("I134").Value = 6086.5713
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 18259.7139
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -15724.7139
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# row 131
$ws.Range("H131").Value = 883.73
$ws.Range("I131").Value = 499
$ws.Range("J131").Value = 926.4778
$ws.Range("K131").Value = 1497
$ws.Range("L131").Value = 2779.4334
$ws.Range("M131").Value = 3543
$ws.Range("N131").Value = -12859.4334

# row 141
$ws.Range("H141").Value = 6788.3335
$ws.Range("I141").Value = 2517.2727
$ws.Range("J141").Value = 13500
$ws.Range("K141").Value = 7551.8181
$ws.Range("L141").Value = 40500
$ws.Range("M141").Value = -2371.8181
$ws.Range("N141").Value = -50860

$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 1638.05
$ws.Range("I97").Value = 1742.3529
$ws.Range("J97").Value = 1047
$ws.Range("K97").Value = 1742.3529
$ws.Range("L97").Value = 1047
$ws.Range("M97").Value = -1246.3529
$ws.Range("N97").Value = -2039

# row 102
$ws.Range("H102").Value = 2058
$ws.Range("I102").Value = 1365.7646
$ws.Range("J102").Value = 3234.8
$ws.Range("K102").Value = 1365.7646
$ws.Range("L102").Value = 3234.8
$ws.Range("M102").Value = 256.2354
$ws.Range("N102").Value = -6478.8

# row 107
$ws.Range("H107").Value = 794.36365
$ws.Range("I107").Value = 840.2
$ws.Range("J107").Value = 756.1667
$ws.Range("K107").Value = 840.2
$ws.Range("L107").Value = 756.1667
$ws.Range("M107").Value = 1079.8
$ws.Range("N107").Value = -4596.1667

# row 113
$ws.Range("H113").Value = 1876.8334
$ws.Range("I113").Value = 1717.375
$ws.Range("K113").Value = 1717.375
$ws.Range("M113").Value = 452.625

$ws = $wb.Worksheets.Item("LTW")
# row 61
$ws.Range("H61").Value = 2985.5715
$ws.Range("I61").Value = 2429.8
$ws.Range("J61").Value = 4375
$ws.Range("K61").Value = 2429.8
$ws.Range("L61").Value = 4375
$ws.Range("M61").Value = -2227.8
$ws.Range("N61").Value = -4779

# row 113
$ws.Range("H113").Value = 2985.5715
$ws.Range("I113").Value = 2429.8
$ws.Range("J113").Value = 4375
$ws.Range("K113").Value = 2429.8
$ws.Range("L113").Value = 4375
$ws.Range("M113").Value = -259.8000000000002
$ws.Range("N113").Value = -8715
